$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.822.91'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.639.46'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.91'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.866.41'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.635.77'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.529'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.72'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.824.89'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '219.17'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.72'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.78%  '
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('E23').Value = '  +3.28%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.16'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.31'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.37'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.63%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.76'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0503'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.259.13'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0178'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.832'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  +2.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.781.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '61.76'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.94%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.10'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '91.88'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('E48').Value = '  +27.11%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.59'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0963'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.26%  '
